$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Swap the "Mobile Number" / "Favorite Color" headers in B1 and C1
$ws.Range("B1").Value = "Favorite Color"
$ws.Range("C1").Value = "Mobile Number"

# Move the active selection from F8 to C2
$ws.Range("C2").Select()
